$wb = $excel.ActiveWorkbook

# ---- Sheet1: Overal Stats -> column AM (new day 43933) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("AL1").Copy()
$ws1.Range("AM1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("AM1").Value = 43933
$ws1.Range("AM3").Value = 10934
$ws1.Range("AM4").Value = 1955
$ws1.Range("AM5").Value = 52
$ws1.Range("AM6").Value = 507
$ws1.Range("AM8").Value = 105
$ws1.Range("AM9").Value = 444
$ws1.Range("AM10").Value = 218
$ws1.Range("AM11").Value = 226
$ws1.Range("AM63").Value = 48
$ws1.Range("AM64").Value = 248
$ws1.Range("AM65").Value = 248
$ws1.Range("AM67").Value = 38
$ws1.Range("AM68").Value = 33
$ws1.Range("AM69").Value = 71
$ws1.Range("AM70").Value = 9
$ws1.Range("AM72").Value = 21
$ws1.Range("AM73").Value = 144
$ws1.Range("AM74").Value = 144
$ws1.Range("AM75").Value = 2

# ---- Sheet2: Total Cases by Ward -> column N (new day 43933) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("M2").Copy()
$ws2.Range("N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("N2").Value = 43933
$ws2.Range("N3").Value = 223
$ws2.Range("N4").Value = 180
$ws2.Range("N5").Value = 146
$ws2.Range("N6").Value = 312
$ws2.Range("N7").Value = 257
$ws2.Range("N8").Value = 288
$ws2.Range("N9").Value = 292
$ws2.Range("N10").Value = 228
$ws2.Range("N11").Value = 29

# ---- Sheet3: Total Cases by Race -> column I (new day 43933) ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("H2").Copy()
$ws3.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws3.Range("I2").Value = 43933
$ws3.Range("I4").Value = 1955
$ws3.Range("I5").Value = 391
$ws3.Range("I6").Value = 365
$ws3.Range("I7").Value = 884
$ws3.Range("I8").Value = 29
$ws3.Range("I9").Value = 6
$ws3.Range("I10").Value = 1
$ws3.Range("I11").Value = 256
$ws3.Range("I12").Value = 23
$ws3.Range("I14").Value = 496
$ws3.Range("I15").Value = 263
$ws3.Range("I16").Value = 1191
$ws3.Range("I17").Value = 5

# ---- Sheet4: Lives Lost by Race -> column I (new day 43933) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("H1").Copy()
$ws4.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws4.Range("I1").Value = 43933
$ws4.Range("I3").Value = 52
$ws4.Range("I4").Value = 2
$ws4.Range("I5").Value = 38
$ws4.Range("I6").Value = 6
$ws4.Range("I7").Value = 6
$ws4.Range("I8").Value = 0

# ---- View updates ----
$ws1.Activate()
$excel.ActiveWindow.Zoom = 110
$ws1.Range("AM63:AM75").Select()

$ws2.Range("N3:N11").Select()

$ws3.Range("I17").Select()

$ws4.Range("I9").Select()

$ws1.Activate()
